$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Source data correction: the "NS" effort entry in row 48 (column C)
# increases from 4 to 5 hours. This is the only underlying data edit;
# the SUMIF/SUM totals and percentage formulas in rows 53 and 56
# recalculate automatically from this single change.
$ws.Range("C48").Value = 5

# Reflect the updated scroll position / active selection captured in the
# saved view (user scrolled down a bit and had moved the selection to I53).
$excel.ActiveWindow.ScrollRow = 31
$excel.ActiveWindow.ScrollColumn = 2
$ws.Range("I53").Select()
